$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.494.57"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.68%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.165.82"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.77%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "164.33"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.14%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.587"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.118"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.65"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.386"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.45%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.716.41"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.65%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.128"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "64.571.49"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.60%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.37"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.66%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.163.54"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.44%  "
$ws.Range("E17").Value = "  -2.45%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "409.48"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.76"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.01%  "
$ws.Range("E20").Value = "  -1.84%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.11"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.22%  "
$ws.Range("E22").Value = "  -0.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.75"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.89%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.485"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.28%  "
$ws.Range("E25").Value = "  -2.33%  "
$ws.Range("E26").Value = "  -6.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.95"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.49%  "
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.82"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.88%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "21.30"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.29%  "
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.93"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.27%  "
$ws.Range("B32").Value = "Aptos"
$ws.Range("C32").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.37"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.00%  "
$ws.Range("E33").Value = "  -1.26%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "155.75"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.33%  "
$ws.Range("E35").Value = "  -2.35%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.696.12"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.71%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.70"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.75%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "24.19"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.10"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.697"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.20%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0623"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.25%  "
$ws.Range("E42").Value = "  -4.98%  "
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "293.01"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.96%  "
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.57"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.91%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0258"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.05%  "
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0987"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.09%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.94"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -8.18%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.49"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.64%  "
$ws.Range("E50").Value = "  -1.75%  "
$ws.Range("E51").Value = "  -5.49%  "
